# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet named "2022-Q1" right before the "总计" sheet,
#    carrying the per-fund breakdown for the new quarter.
# 2) Insert a new top data row in "总计" summarising that quarter
#    (date / holding count / holding market value).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New "2022-Q1" sheet, inserted just before "总计"
# ---------------------------------------------------------------------
$totalSheetRef = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheetRef)
$newSheet.Name = "2022-Q1"

# Borrow formatting (header row style + column-A index style) from an
# existing quarter sheet so the new sheet matches the rest of the workbook.
$template = $wb.Worksheets.Item("2021-Q4")
$template.Range("A1:H1").Copy()
$newSheet.Range("A1:H1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$newSheet.Range("A2:A10").PasteSpecial(-4122)

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Make sure the fund code column and the numeric-looking fund figures are
# stored as text (to match the source data, which keeps these as
# formatted strings) while the index column (A) and rank column (H) stay
# numeric. Must be set *before* the values are assigned, otherwise
# digit-only strings (e.g. fund code "002415") get coerced to numbers and
# lose their leading zeros.
$newSheet.Range("B2:G10").NumberFormat = "@"

$fundRows = @(
    @("002415", "融通通盈灵活配置混合",           "0.89", "69.03", "8.01", "0.0713", 1),
    @("515450", "南方标普中国A股大盘红利低波50ETF", "2.04", "99.51", "2.83", "0.0577", 8),
    @("008114", "天弘中证红利低波动100指数A",       "3.16", "92.60", "1.71", "0.0540", 6),
    @("010746", "富安达长三角区域主题混合",         "1.16", "92.45", "4.49", "0.0521", 10),
    @("008115", "天弘中证红利低波动100指数C",       "2.37", "92.60", "1.71", "0.0405", 6),
    @("008091", "中信保诚红利精选混合A",           "0.90", "90.17", "3.99", "0.0359", 3),
    @("515100", "景顺长城中证红利低波动100ETF",     "1.25", "97.96", "1.82", "0.0228", 6),
    @("008092", "中信保诚红利精选混合C",           "0.55", "90.17", "3.99", "0.0219", 3),
    @("162205", "泰达宏利风险预算混合",             "1.34", "27.82", "0.70", "0.0094", 10)
)

$row = 2
foreach ($fund in $fundRows) {
    $newSheet.Range("A$row").Value = $row - 2
    $newSheet.Range("B$row").Value = $fund[0]
    $newSheet.Range("C$row").Value = $fund[1]
    $newSheet.Range("D$row").Value = $fund[2]
    $newSheet.Range("E$row").Value = $fund[3]
    $newSheet.Range("F$row").Value = $fund[4]
    $newSheet.Range("G$row").Value = $fund[5]
    $newSheet.Range("H$row").Value = $fund[6]
    $row = $row + 1
}

# ---------------------------------------------------------------------
# 2) Add the 2022-Q1 summary row to "总计" (goes right under the header,
#    above the existing 2021-Q4 row — rows stay in reverse-chronological
#    order). Re-fetch the sheet by name: after inserting the new sheet
#    above, any previously-held reference to "总计" now points at the
#    newly inserted sheet instead.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# The existing rows' "A" column is a simple 0-based row counter; bump each
# one by 1 to make room for the new row 2 (done before inserting, counting
# up from the bottom so the writes don't clobber each other).
for ($r = 6; $r -ge 2; $r--) {
    $oldValue = $totalSheet.Cells.Item($r, 1).Value2
    $totalSheet.Cells.Item($r, 1).Value = $oldValue + 1
}

$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 9
$totalSheet.Range("D2").Value = 0.37
